# regen sval data to filter save games
# Updates the numeric stat columns (B:E, G) on Sheet1 with newly
# regenerated values. Column A (date) and column F (win flag) are
# unchanged; column G is the row sum of B:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    3  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    4  = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    5  = @{ B = 0.2881169905109251; C = 208501.5462402375;   D = 0.7210945179870265; E = 13.86384647080068;  G = 208516.4192982168 }
    6  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    7  = @{ B = 0.6545652718822623; C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987; G = 3.536033448013082 }
    8  = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    9  = @{ B = 0.6545652718822623; C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987; G = 3.536033448013082 }
    10 = @{ B = 1.445647641019636;  C = 0.3048912486333797;  D = 0.7210945179870265; E = 13.86384647080068;  G = 16.33547987844073 }
    11 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 18.71679738969934;  E = 0.5333859586016987; G = 24.14949828602258 }
    12 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 13.86384647080068;  G = 19.48425592650926 }
    13 = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 18.71679738969934;  E = 0.5333859586016987; G = 22.32281868886277 }
    14 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in @("B", "C", "D", "E", "G")) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
